$d = $word.ActiveDocument

# "Login: se cambia el login para que funcione con el correo electrónico."
# The title under the project is currently "Instructora:" (feminine) and
# must be changed to "Instructor:" (masculine) to reflect the new
# instructor's name shown right below it.
$d.Content.Find.Execute("Instructora:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Instructor:", 2)
